$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a daily price log: newest date lives in row 2, older dates
# further down. Today's pull adds a new "10-11-2025" row at the top and
# pushes every existing data row (2..96) down by one (-> 3..97); the row
# that falls off the bottom (old row 96, 07-08-2025) reappears duplicated
# at the new row 97, mirroring the upstream feed's own duplicate entry.

$lastRow = 96
$newLastRow = 97

# 1) Duplicate the last row's formatting (borders/alignment/number format)
#    into the brand-new row so it matches the rest of the table instead of
#    falling back to unstyled defaults.
$ws.Range("A$lastRow`:F$lastRow").Copy($ws.Range("A$newLastRow`:F$newLastRow"))

# 2) Snapshot the current (pre-shift) contents of rows 2..96 as plain text
#    /numbers before overwriting anything, so later writes can't clobber a
#    value we still need to read.
$dates = @{}
$descs = @{}
$codes = @{}
$prices = @{}
$circDates = @{}
$links = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $dates[$r]     = $ws.Cells.Item($r, 1).Text
    $descs[$r]     = $ws.Cells.Item($r, 2).Text
    $codes[$r]     = $ws.Cells.Item($r, 3).Text
    $prices[$r]    = $ws.Cells.Item($r, 4).Value2
    $circDates[$r] = $ws.Cells.Item($r, 5).Text
    $links[$r]     = $ws.Cells.Item($r, 6).Text
}

# 3) Write the shifted data into rows 3..97, bottom-to-top so a source row
#    is always read before it gets overwritten. Date-like text (A, E) is
#    apostrophe-prefixed on entry so COM doesn't silently reinterpret it
#    as a real date serial; the prefix is stripped again in step 6.
for ($r = $lastRow; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value = "'" + $dates[$r]
    $ws.Cells.Item($dest, 2).Value = $descs[$r]
    $ws.Cells.Item($dest, 3).Value = $codes[$r]
    $ws.Cells.Item($dest, 4).Value2 = $prices[$r]
    $ws.Cells.Item($dest, 5).Value = "'" + $circDates[$r]
    $ws.Cells.Item($dest, 6).Value = $links[$r]
}

# 4) Write the brand-new top row (today's price pull). Same circular/price
#    /link as the old row 2 since no new circular was issued today.
$ws.Cells.Item(2, 1).Value = "'10-11-2025"
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value2 = 297.15
$ws.Cells.Item(2, 5).Value = "'01-11-2025"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

# 5) Hyperlinks don't follow cell shifts automatically in this host, and a
#    per-cell Hyperlinks.Delete() wipes the *entire* collection, so rebuild
#    the full F2:F$newLastRow hyperlink set from the (now-correct) link
#    text in one pass.
$ws.Cells.Hyperlinks.Delete()
for ($r = 2; $r -le $newLastRow; $r++) {
    $target = $ws.Cells.Item($r, 6).Text
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
}

# 6) Cosmetic cleanup: the apostrophe-prefix trick and Hyperlinks.Add both
#    nudge a cell's style (quote-prefix flag / auto blue-underline font).
#    Column C always keeps the sheet's plain centered style, so copy its
#    formatting (values untouched) back onto A, E and F.
$ws.Range("C2:C$newLastRow").Copy()
$ws.Range("A2:A$newLastRow").PasteSpecial(-4122)
$ws.Range("E2:E$newLastRow").PasteSpecial(-4122)
$ws.Range("F2:F$newLastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
